$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 39
$ws.Range("F5").Value = 440
$ws.Range("F6").Value = 1867
$ws.Range("F8").Value = 1352
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = 756
$ws.Range("F11").Value = 128
$ws.Range("F12").Value = 2906
$ws.Range("F13").Value = 385
$ws.Range("F14").Value = 879
$ws.Range("F15").Value = 1130
$ws.Range("F16").Value = 591
$ws.Range("F17").Value = 334
$ws.Range("F18").Value = 69
$ws.Range("F19").Value = 1660
$ws.Range("F20").Value = 343
$ws.Range("F21").Value = 1264
$ws.Range("F22").Value = 212
$ws.Range("F23").Value = 600
$ws.Range("F25").Value = 1072
$ws.Range("F26").Value = 1520
$ws.Range("F27").Value = 1472
$ws.Range("F28").Value = 1340
$ws.Range("F29").Value = 339
$ws.Range("F30").Value = 1293
$ws.Range("F31").Value = 445
$ws.Range("F32").Value = 153
$ws.Range("F33").Value = 974
$ws.Range("F35").Value = 1855
$ws.Range("F36").Value = 483
$ws.Range("F37").Value = 1053
$ws.Range("F38").Value = 157
$ws.Range("F40").Value = 2295
$ws.Range("F41").Value = 151
$ws.Range("F42").Value = 892
$ws.Range("F43").Value = 2800
$ws.Range("F48").Value = 10

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 62
$ws.Range("F10").Value = 36
$ws.Range("F12").Value = 369
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 47
$ws.Range("F16").Value = 21
$ws.Range("F17").Value = 70
$ws.Range("F18").Value = 70
$ws.Range("F20").Value = 288
$ws.Range("F23").Value = 67
$ws.Range("F25").Value = 70
$ws.Range("F26").Value = 63
$ws.Range("F27").Value = 63
$ws.Range("F29").Value = 9
$ws.Range("F31").Value = 128
$ws.Range("F38").Value = 162

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 295
$ws.Range("F5").Value = 3029
$ws.Range("F6").Value = 4857
$ws.Range("F9").Value = 678
$ws.Range("F10").Value = 948
$ws.Range("F11").Value = 550
$ws.Range("F12").Value = 640
$ws.Range("F13").Value = 1368
$ws.Range("F15").Value = 1274

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 295
$ws.Range("F4").Value = 4857
$ws.Range("F5").Value = 678
$ws.Range("F6").Value = 948
$ws.Range("F7").Value = 550
$ws.Range("F9").Value = 640
$ws.Range("F10").Value = 1368
$ws.Range("F11").Value = 1867
$ws.Range("F13").Value = 1352
$ws.Range("F14").Value = 756
$ws.Range("F15").Value = 756
$ws.Range("F16").Value = 1274
$ws.Range("F17").Value = 2906
$ws.Range("F18").Value = 36
$ws.Range("F19").Value = 385
$ws.Range("F20").Value = 879
$ws.Range("F21").Value = 1130
$ws.Range("F22").Value = 591
$ws.Range("F23").Value = 334
$ws.Range("F24").Value = 1660
$ws.Range("F26").Value = 343
$ws.Range("F27").Value = 369
$ws.Range("F28").Value = 1264
$ws.Range("F29").Value = 212
$ws.Range("F30").Value = 600
$ws.Range("F32").Value = 1520
$ws.Range("F33").Value = 1472
$ws.Range("F34").Value = 1340
$ws.Range("F35").Value = 339
$ws.Range("F36").Value = 70
$ws.Range("F37").Value = 1293
$ws.Range("F38").Value = 445
$ws.Range("F39").Value = 974
$ws.Range("F41").Value = 1855
$ws.Range("F42").Value = 63
$ws.Range("F43").Value = 128
$ws.Range("F45").Value = 2295
$ws.Range("F46").Value = 892
$ws.Range("F47").Value = 2800
